$d = $word.ActiveDocument

# wdAlignParagraphJustify = 3
$wdAlignParagraphJustify = 3

# --- 1. Justify the paragraphs that are untouched structurally (indices stable) ---
# Paragraph 4: "En esta ocasión..."
$d.Paragraphs(4).Alignment = $wdAlignParagraphJustify
# Paragraph 5: "Muestra de esto..."
$d.Paragraphs(5).Alignment = $wdAlignParagraphJustify
# Paragraph 6: "El ingreso a las instalaciones..."
$d.Paragraphs(6).Alignment = $wdAlignParagraphJustify
# Paragraph 7: "Posteriormente se daba..."
$d.Paragraphs(7).Alignment = $wdAlignParagraphJustify
# Paragraph 8: "Otra área de mucha participación..."
$d.Paragraphs(8).Alignment = $wdAlignParagraphJustify
# Paragraph 9: "Nos encontramos con empresas..."
$d.Paragraphs(9).Alignment = $wdAlignParagraphJustify
# Paragraph 10: empty paragraph
$d.Paragraphs(10).Alignment = $wdAlignParagraphJustify

# --- 2. Merge paragraph 11 ("...Epson, ") with paragraph 12 ("Robot limpiador...") ---
# by deleting paragraph 11's trailing paragraph mark. (The merged paragraph
# inherits paragraph 12's pPr, so alignment must be (re)applied afterwards.)
$p11 = $d.Paragraphs(11)
$markRange = $d.Range($p11.Range.End - 1, $p11.Range.End)
$markRange.Delete()
$d.Paragraphs(11).Alignment = $wdAlignParagraphJustify

# --- 3. Replace the old "Robot limpiador..." sentence (now tail of paragraph 11)
#        with the new closing sentence for that paragraph.
$d.Content.Find.Execute(
    "Robot limpiador de casa y muestra de robots ganadores en torneos de robotica.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "y otras más, pero no estaban los representantes directos de la marca, sino los distribuidores o vendedores al usuario final, lo que hace que los artículos no tengan precio directo.",
    2
) | Out-Null

# --- 4. The paragraph that used to hold only the _GoBack bookmark is now paragraph 12. ---
$p12 = $d.Paragraphs(12)
$p12.Alignment = $wdAlignParagraphJustify

# --- 5. Insert the new narrative text around the bookmark, preserving its position. ---
#        NOTE: the "after" insertion must happen before the "before" insertion,
#        otherwise the zero-width bookmark's anchor drifts and a later
#        InsertAfter lands on the wrong side.
$bm = $d.Bookmarks("_GoBack")
$bmr = $bm.Range
$bmr.InsertAfter("IPN y el Tecnológico de Monterrey.")

$bm2 = $d.Bookmarks("_GoBack")
$bmr2 = $bm2.Range
$bmr2.InsertBefore("Se presentó también como novedad el Robot Aspiradora, que realiza la colecta de polvo y pelusa en casa de forma autónoma, además se realizó una exhibición de los robots ganadores en torneos de robótica en los que asistió el ")

# --- 6. Add a new empty paragraph at the very end of the document. ---
#        (It must NOT inherit the "both" justification of the paragraph
#        before it, so reset its alignment back to the default afterwards.)
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$wdAlignParagraphLeft = 0
$newLastPara = $d.Paragraphs($d.Paragraphs.Count)
$newLastPara.Alignment = $wdAlignParagraphLeft

Write-Output "done"
